$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.879.72'
$ws.Range("E2").Value = '  -1.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.358.03'
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.30'
$ws.Range("E5").Value = '  +1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.27'
$ws.Range("E6").Value = '  -4.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.642'
$ws.Range("E7").Value = '  +0.43%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.625'
$ws.Range("E9").Value = '  -1.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.23'

$ws.Range("E11").Value = '  -1.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.53'
$ws.Range("E12").Value = '  -2.11%  '

$ws.Range("E13").Value = '  -3.58%  '

$ws.Range("E14").Value = '  +0.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.19'
$ws.Range("E15").Value = '  -2.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.714.94'
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.365.73'
$ws.Range("E17").Value = '  -2.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.836.46'
$ws.Range("E18").Value = '  -1.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.92'

$ws.Range("E20").Value = '  -2.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '77.07'
$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.65'
$ws.Range("E22").Value = '  +4.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '265.48'
$ws.Range("E23").Value = '  +0.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.33'
$ws.Range("E24").Value = '  -7.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.01'
$ws.Range("E25").Value = '  +9.05%  '

$ws.Range("E26").Value = '  +0.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.53'
$ws.Range("E27").Value = '  -4.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.06'
$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("E29").Value = '  -1.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.81'
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.14'
$ws.Range("E31").Value = '  -2.09%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.26'
$ws.Range("E32").Value = '  +4.96%  '

$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.46'
$ws.Range("E34").Value = '  -9.08%  '

$ws.Range("E35").Value = '  +1.70%  '

$ws.Range("E36").Value = '  +7.04%  '

$ws.Range("E37").Value = '  -7.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0361'
$ws.Range("E38").Value = '  -3.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.80'
$ws.Range("E39").Value = '  -7.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.73'
$ws.Range("E40").Value = '  -3.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.238'
$ws.Range("E41").Value = '  +2.59%  '

$ws.Range("E42").Value = '  -1.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.85'
$ws.Range("E43").Value = '  -1.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.48'
$ws.Range("E44").Value = '  +25.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '120.70'
$ws.Range("E45").Value = '  +6.98%  '

$ws.Range("E46").Value = '  -0.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.96'
$ws.Range("E47").Value = '  -5.11%  '

$ws.Range("E48").Value = '  -0.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.17'
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("E50").Value = '  -3.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.101'
$ws.Range("E51").Value = '  +0.01%  '
